$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sections")

$ws.Range("AM2:AN3").ClearContents()
